$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7877166271209717
$ws.Range("B1").Value = 1.809979915618896
$ws.Range("C1").Value = 2.545440196990967
$ws.Range("D1").Value = 1.528391599655151
$ws.Range("E1").Value = 0.8194904923439026
